$d = $word.ActiveDocument

# Locate (by index) the paragraph that currently reads
# "4. Execute the procedure on a target machine."
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("4. Execute the procedure")) {
        $targetIndex = $i
        break
    }
}

# Insert a brand-new empty paragraph immediately before the target paragraph; this
# will become the new step "4" (the optional LogIt-variable instructions). The
# target paragraph is pushed down by one index and will become step "5".
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()

$newIndex = $targetIndex
$shiftedTargetIndex = $targetIndex + 1

# Fill the newly created (now-empty) paragraph with the rich content, using
# InsertXML so per-run character formatting (italic / bold) and the
# spell-check-exception markers around "LogIt" can be set in a single shot.
$newRange = $d.Paragraphs.Item($newIndex).Range
$newParaXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">4. (Optional) </w:t></w:r>
            <w:r><w:t>set</w:t></w:r>
            <w:r><w:t xml:space="preserve"> the </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>LogIt</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> variable </w:t></w:r>
            <w:r><w:t xml:space="preserve">to 1 if you need extended logging </w:t></w:r>
            <w:r><w:t xml:space="preserve">information </w:t></w:r>
            <w:r><w:t>for the Procedure</w:t></w:r>
            <w:r><w:t xml:space="preserve"> in the </w:t></w:r>
            <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Procedure History</w:t></w:r>
            <w:r><w:t>. Otherwise set it to 0.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$newRange.InsertXML($newParaXml)

# Finally, renumber the (still intact) target paragraph from "4" to "5" by swapping
# just its leading digit, leaving the remaining text untouched.
$shiftedRange = $d.Paragraphs.Item($shiftedTargetIndex).Range
$digitStart = $shiftedRange.Start
$digitRange = $d.Range($digitStart, $digitStart + 1)
$digitRange.Delete()
$insertionPoint = $d.Range($digitStart, $digitStart)
$insertionPoint.InsertAfter("5")

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i => [$($d.Paragraphs.Item($i).Range.Text)]"
}
